$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the placeholder string "0pfj7wGzbgUD" with "XXXXXX" across column C (rows 2-95)
$ws.Range("C2:C95").Value = "XXXXXX"

# Update the sheet view / selection to match the committed state
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()
